# Fill in the missing CENTRODECOSTO (column F) values for rows 40-81.
# These cells were previously blank and are now populated with the
# appropriate cost-center text, matching the same values used for the
# corresponding DOCUMENTO/ORGANIZACION rows elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    40 = "Administración"
    41 = "Depósito"
    42 = "Depósito"
    43 = "Administración"
    44 = "Administración"
    45 = "Depósito"
    46 = "Ventas"
    47 = "Depósito"
    48 = "Depósito"
    49 = "Ventas"
    50 = "Depósito"
    51 = "Ventas"
    52 = "Depósito"
    53 = "Administración"
    54 = "Administración"
    55 = "Depósito"
    56 = "Administración"
    57 = "Ventas"
    58 = "Ventas"
    59 = "Depósito"
    60 = "Administración"
    61 = "Depósito"
    62 = "Depósito"
    63 = "Depósito"
    64 = "Depósito"
    65 = "Depósito"
    66 = "Administración"
    67 = "Depósito"
    68 = "Administración"
    69 = "Ventas"
    70 = "Ventas"
    71 = "Administración"
    72 = "Ventas"
    73 = "Ventas"
    74 = "Depósito"
    75 = "Depósito"
    76 = "Depósito"
    77 = "Administración"
    78 = "Depósito"
    79 = "Ventas"
    80 = "Administración"
    81 = "Administración"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $values[$row]
    $cell.NumberFormat = "@"
}
